# Saldo.xlsx — "Add files via upload"
#
# The uploaded sheet reshuffled one data row: the row for account
# 005186167 / ANDREA used to sit right after 005003629 / ANDRE (row 21)
# with a balance of 1674.77. In the new upload that row was removed from
# there and re-appended at the end of the same balance-ordered block,
# just after 005646524 / EVANGELINA, with its balance corrected to
# 674.77 (looks like a stray leading "1" was dropped).
#
# Reproduce that with a plain cut/paste: copy the whole row to its new
# location (so the text-formatted account/name cells keep their exact
# formatting instead of Excel re-guessing the type of the pasted text),
# delete it from its old location, then fix up the corrected balance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the ANDREA / 005186167 row (row 21) and insert that copy right
#    after 005646524 / EVANGELINA (row 31), i.e. as the new row 32.
$ws.Rows(21).Copy()
$ws.Rows(32).Insert()

# 2) Remove the now-duplicated row from its original spot (row 21).
#    Everything below it — including the pasted copy — shifts up one, so
#    the copy ends up at row 31, right after EVANGELINA.
$ws.Rows(21).Delete()

# 3) Correct the balance on the relocated row.
$ws.Cells.Item(31, 3).Value = 674.77
